$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 89.6727272727273
$ws.Range("C4").Value = 0.894604392662199

$ws.Range("B7").Value = 96.8727272727273
$ws.Range("C7").Value = 0.968066725825364

$ws.Range("B8").Value = 96.5454545454545
$ws.Range("C8").Value = 0.964745911470338

$ws.Range("C10").Value = 0.965107477069471

$ws.Range("C11").Value = 0.960284124069562

$ws.Range("B16").Value = 93.6
$ws.Range("C16").Value = 0.934709505864273

$ws.Range("B17").Value = 91.9272727272727
$ws.Range("C17").Value = 0.917628299489201
